# Generate Report for Archive
# Update the localization status from "Ready for handoff" to "In Translation"
# across the Overview / zh-cn / de-de sheets, then tighten the Status
# column widths to match the shorter text (the report generator re-ran its
# auto-fit pass over the narrower "In Translation" values).

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"
# Target stored width is ~13.41 "characters"; the column-width grid here
# snaps to 1/6 increments, so 12.5 is the closest achievable input.
$statusColWidth = 12.5

# --- Overview sheet: zh-cn (E) and de-de (F) status columns ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsOverview.Range("E2:F3").ColumnWidth = $statusColWidth

# --- zh-cn sheet: Status column (C) ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus
$wsZh.Range("C2:C3").ColumnWidth = $statusColWidth

# --- de-de sheet: Status column (C) ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus
$wsDe.Range("C2:C3").ColumnWidth = $statusColWidth
